# Annotate OneSpin FV Results
# Populate the plans w/ tests and results
#
# This script updates the CV32E40P_RV32M_Extension_Instructions workbook:
#  - RV32M sheet: "Properties" (col F), "Coverage Method" (col I), "Coverage
#    Details" (col J) and "Link to Coverage" (col K) are populated with the
#    actual OneSpin formal-verification check names / results instead of the
#    "TBD" placeholders.
#  - DONOTDELETE sheet: the drop-down source list in column A is updated from
#    the placeholder RV32I/RVC check names to the real RV_chk.ops.RV32M.*
#    check names (plus RV32M.all_a), and the now-unused trailing rows are
#    removed.

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("RV32M")
$wsList = $wb.Worksheets.Item("DONOTDELETE")

# ---------------------------------------------------------------------------
# 1) RV32M sheet - populate Properties / Coverage Method / Coverage Details /
#    Link to Coverage for every instruction row (2-9).
# ---------------------------------------------------------------------------

$longDetails = "Each operand is controlled by a parameter defining the number of bits to be tied to '0 or '1:" + [char]10 + "(opA[XLEN-1:PARAM]='0 || opA[XLEN-1:PARAM]='1) && (opB[XLEN-1:PARAM]='0 || opB[XLEN-1:PARAM]='1) "
$shortDetails = "Each operand is controlled by a parameter defining the number of bits to be tied to '0 or '1 (default is 1)"

# row -> Properties (col F) value
$propMap = @{
    2 = "RV_chk.ops.RV32M.mul_a"
    3 = "RV_chk.ops.RV32M.mulh_a"
    4 = "RV_chk.ops.RV32M.mulhu_a"
    5 = "RV_chk.ops.RV32M.mulhsu_a"
    6 = "RV_chk.ops.RV32M.div_a"
    7 = "RV_chk.ops.RV32M.rem_a"
    8 = "RV_chk.ops.RV32M.divu_a"
    9 = "RV_chk.ops.RV32M.remu_a"
}

foreach ($row in 2..9) {
    # Properties (F): TBD -> RV_chk.ops.RV32M.<op>_a
    $wsData.Cells.Item($row, 6).Value = $propMap[$row]

    # Coverage Method (I): Assertion -> Partial Proof
    $wsData.Cells.Item($row, 9).Value = "Partial Proof"

    # Coverage Details (J): placeholder -> real coverage-details text
    if ($row -eq 2) {
        $wsData.Cells.Item($row, 10).Value = $longDetails
    } else {
        $wsData.Cells.Item($row, 10).Value = $shortDetails
    }
    $wsData.Cells.Item($row, 10).WrapText = $true

    # Link to Coverage (K): TBD -> N/A
    $wsData.Cells.Item($row, 11).Value = "N/A"
}

# Row heights settle slightly differently once the text/column widths change.
$wsData.Rows.Item(2).RowHeight = 55.2
$wsData.Rows.Item(3).RowHeight = 44.4
$wsData.Rows.Item(4).RowHeight = 42

# Column widths: Properties (F) and Link to Coverage (G... wait see below) and
# Coverage Details (J) need to grow to fit the new text.
$wsData.Columns.Item(6).ColumnWidth = 30.5546875
$wsData.Columns.Item(7).ColumnWidth = 19.33203125
$wsData.Columns.Item(10).ColumnWidth = 66.109375

# Update the view: scroll so column D is left-most and select G7 (matches the
# state the workbook was saved in).
$wsData.Activate()
$wsData.Range("G7").Select()

# ---------------------------------------------------------------------------
# 2) DONOTDELETE sheet - refresh the drop-down source list (column A) with the
#    real RV_chk.ops.RV32M.* check names, keep RV32M.all_a, and drop the
#    leftover RVC.* placeholder rows (the list shrinks from 15 to 11 entries).
# ---------------------------------------------------------------------------

$propList = @(
    "RV_chk.ops.RV32M.div_a",
    "RV_chk.ops.RV32M.divu_a",
    "RV_chk.ops.RV32M.mul_a",
    "RV_chk.ops.RV32M.mulh_a",
    "RV_chk.ops.RV32M.mulhsu_a",
    "RV_chk.ops.RV32M.mulhu_a",
    "RV_chk.ops.RV32M.rem_a",
    "RV_chk.ops.RV32M.remu_a",
    "RV32M.all_a"
)

foreach ($i in 0..($propList.Length - 1)) {
    $wsList.Cells.Item($i + 4, 1).Value = $propList[$i]
}

# Remove the now-unused rows 13-16 (RVC.Arith_a .. RVC.Mem_a)
$wsList.Range("A13:G16").EntireRow.Delete()

$wsList.Columns.Item(1).ColumnWidth = 29.44140625
$wsList.Range("A4").Select()
